$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 8999.75
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 11666.333
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 11666.333
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -12318.333

# Row 33
$ws.Range("H33").Value = 349
$ws.Range("I33").Value = 304.44446
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 304.44446
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -75.44445999999999
$ws.Range("N33").Value = -1208

# Row 74
$ws.Range("H74").Value = 3499.75
$ws.Range("I74").Value = 2999.6667
$ws.Range("K74").Value = 2999.6667
$ws.Range("M74").Value = -2063.6667

# Row 77
$ws.Range("H77").Value = 3499.75
$ws.Range("I77").Value = 2999.6667
$ws.Range("K77").Value = 14998.3335
$ws.Range("M77").Value = -10318.3335

# Row 92
$ws.Range("H92").Value = 1253.1111
$ws.Range("I92").Value = 1002.8
$ws.Range("K92").Value = 1002.8
$ws.Range("M92").Value = 245.2

# Row 100
$ws.Range("H100").Value = 9331.333000000001
$ws.Range("I100").Value = 2982
$ws.Range("J100").Value = 10125
$ws.Range("K100").Value = 2982
$ws.Range("L100").Value = 10125
$ws.Range("M100").Value = -2441
$ws.Range("N100").Value = -11207

# Row 132
$ws.Range("H132").Value = 3238.375
$ws.Range("I132").Value = 2272.5715
$ws.Range("K132").Value = 6817.7145
$ws.Range("M132").Value = -4287.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16541.273
$ws.Range("I32").Value = 16541.273
$ws.Range("K32").Value = 16541.273
$ws.Range("M32").Value = -16254.273

# Row 80
$ws.Range("H80").Value = 78110
$ws.Range("J80").Value = 78110
$ws.Range("L80").Value = 78110
$ws.Range("N80").Value = -80106

# Row 83
$ws.Range("H83").Value = 78110
$ws.Range("J83").Value = 78110
$ws.Range("L83").Value = 234330
$ws.Range("N83").Value = -244314

# Row 137
$ws.Range("H137").Value = 69696
$ws.Range("J137").Value = 69696
$ws.Range("L137").Value = 69696
$ws.Range("N137").Value = -79896

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1249.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2590.2778
$ws.Range("I31").Value = 2607.353
$ws.Range("K31").Value = 2607.353
$ws.Range("M31").Value = -2312.353

# Row 34
$ws.Range("H34").Value = 2590.2778
$ws.Range("I34").Value = 2607.353
$ws.Range("K34").Value = 2607.353
$ws.Range("M34").Value = -2405.353

# Row 58
$ws.Range("H58").Value = 2500
$ws.Range("I58").Value = 2500
$ws.Range("K58").Value = 2500
$ws.Range("M58").Value = -2297

# Row 107
$ws.Range("H107").Value = 1485.6111
$ws.Range("I107").Value = 962
$ws.Range("K107").Value = 962
$ws.Range("M107").Value = 958

# Row 132
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20057

# Row 136
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

# Row 141
$ws.Range("H141").Value = 86085.664
$ws.Range("J141").Value = 84820.73
$ws.Range("L141").Value = 84820.73
$ws.Range("N141").Value = -95180.73

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 225
$ws.Range("J86").Value = 250
$ws.Range("L86").Value = 750
$ws.Range("N86").Value = -3122

# Row 89
$ws.Range("H89").Value = 225
$ws.Range("J89").Value = 250
$ws.Range("L89").Value = 2250
$ws.Range("N89").Value = -14106

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# Row 140
$ws.Range("H140").Value = 1766
$ws.Range("I140").Value = 1618.2858
$ws.Range("K140").Value = 4854.857400000001
$ws.Range("M140").Value = 325.1425999999992

$ws = $wb.Worksheets.Item("GSM")
# Row 139
$ws.Range("H139").Value = 79663.39999999999
$ws.Range("J139").Value = 79663.39999999999
$ws.Range("L139").Value = 79663.39999999999
$ws.Range("N139").Value = -89943.39999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1825.3334
$ws.Range("I7").Value = 1744.5
$ws.Range("K7").Value = 1744.5
$ws.Range("M7").Value = -1632.5

# Row 55
$ws.Range("H55").Value = 448.875
$ws.Range("I55").Value = 424.6
$ws.Range("J55").Value = 489.33334
$ws.Range("K55").Value = 424.6
$ws.Range("L55").Value = 489.33334
$ws.Range("M55").Value = -251.6
$ws.Range("N55").Value = -835.33334

# Row 61
$ws.Range("H61").Value = 3246.9285
$ws.Range("I61").Value = 3150.5386
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 3150.5386
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -2948.5386
$ws.Range("N61").Value = -4904

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 93
$ws.Range("H93").Value = 2672.5
$ws.Range("I93").Value = 2720.7144
$ws.Range("K93").Value = 2720.7144
$ws.Range("M93").Value = -1472.7144

# Row 113
$ws.Range("H113").Value = 3246.9285
$ws.Range("I113").Value = 3150.5386
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 3150.5386
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -980.5385999999999
$ws.Range("N113").Value = -8840

# Row 126
$ws.Range("H126").Value = 1825.3334
$ws.Range("I126").Value = 1744.5
$ws.Range("K126").Value = 5233.5
$ws.Range("M126").Value = -2763.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2360.9412
$ws.Range("I126").Value = 2360.9412
$ws.Range("K126").Value = 7082.823600000001
$ws.Range("M126").Value = -4612.823600000001

# Row 132
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20057

# Row 133
$ws.Range("H133").Value = 29999.5
$ws.Range("J133").Value = 29999.5
$ws.Range("L133").Value = 29999.5
$ws.Range("N133").Value = -40119.5
